$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old rows 2-41 range in column A (rows below 13 will be removed afterward)

$ws.Range("A2").Value = '("Ajani''s Pridemate", [''Token Creature — Cat Soldier'', ''Whenever you gain life, put a +1/+1 counter on Ajani’s Pridemate.'', ''2/2''])'
$ws.Range("A3").Value = '(''Chandra, Awakened Inferno Emblem'', [''Emblem'', ''At the beginning of your upkeep, this emblem deals 1 damage to you.''])'
$ws.Range("A4").Value = '(''Demon'', [''Token Creature — Demon'', ''Flying'', ''5/5''])'
$ws.Range("A5").Value = '(''Elemental'', [''Token Creature — Elemental'', ''1/1''])'
$ws.Range("A6").Value = '(''Elemental Bird'', [''Token Creature — Elemental Bird'', ''Flying'', ''4/4''])'
$ws.Range("A7").Value = '(''Golem'', [''Token Artifact Creature — Golem'', ''3/3''])'
$ws.Range("A8").Value = '(''Mu Yanling, Sky Dancer Emblem'', [''Emblem'', ''Islands you control have “{T}: Draw a card.”''])'
$ws.Range("A9").Value = '(''Soldier'', [''Token Creature — Soldier'', ''1/1''])'
$ws.Range("A10").Value = '(''Spirit'', [''Token Creature — Spirit'', ''Flying'', ''1/1''])'
$ws.Range("A11").Value = '(''Treasure'', [''Token Artifact — Treasure'', ''{T}, Sacrifice this artifact: Add one mana of any color.''])'
$ws.Range("A12").Value = '(''Wolf'', [''Token Creature — Wolf'', ''2/2''])'
$ws.Range("A13").Value = '(''Zombie'', [''Token Creature — Zombie'', ''2/2''])'

# Remove leftover rows 14-41 which no longer have data
$ws.Range("A14:A41").Clear()
